# BL-197: Removed b2b question and updated VIDs covered by automation.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update / add the "Variation" (column E) VID values ---------------
# b2 sheet rows
$ws.Range("E2").Value = "0-0-0-0-1-2-0"
$ws.Range("E3").Value = "0-0-0-0-1-2-0"
$ws.Range("E4").Value = "0-0-0-0-1-2-0"
$ws.Range("E5").Value = "0-0-0-0-1-2-0"
$ws.Range("E6").Value = "0-0-0-0-0-2-0"
$ws.Range("E6").Style = "Normal"
$ws.Range("E7").Value = "0-0-0-0-0-2-0"
$ws.Range("E7").Style = "Normal"

# bl sheet rows
$ws.Range("E8").Value  = "0-0-5-0-0-1-2-0-0"
$ws.Range("E9").Value  = "0-0-5-0-0-1-2-0-0"
$ws.Range("E10").Value = "0-0-5-0-0-1-2-0-0"
$ws.Range("E11").Value = "0-0-5-0-0-1-2-0-0"
$ws.Range("E12").Value = "0-0-5-0-0-0-2-0-0"
$ws.Range("E12").Style = "Normal"
$ws.Range("E13").Value = "0-0-5-0-0-0-2-0-0"
$ws.Range("E13").Style = "Normal"

# --- Swap the BrowserType values between the two Prod bl rows ---------
$ws.Range("C12").Value = "IE"
$ws.Range("C13").Value = "FIREFOX"

# --- Reset the saved selection back to the top of the sheet -----------
$ws.Range("A1").Select()
